# Fix for "computing relative excess deaths": the Percent_excess/percent_cri
# columns (E:F) had been populated with the Rate_excess/rate_cri values (and
# vice versa) - i.e. columns E:F and G:H were swapped. Swap the two column
# pairs back for every row (header + all data rows) in the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: does this text look like something Excel would auto-convert to a
# number/percentage (and so needs a leading quote to stay literal text, just
# like retyping it by hand would)?
function Test-LooksNumeric($s) {
    if ($null -eq $s) { return $false }
    if ($s -eq "") { return $false }
    $clean = $s -replace ',', ''
    if ($clean -match '^-?\d+(\.\d+)?%?$') { return $true }
    return $false
}

function Set-CellText($cell, $text) {
    if ($null -eq $text) {
        $cell.Value = ""
    } elseif (Test-LooksNumeric $text) {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)   # column E - Percent_excess
    $fCell = $ws.Cells.Item($r, 6)   # column F - percent_cri
    $gCell = $ws.Cells.Item($r, 7)   # column G - Rate_excess
    $hCell = $ws.Cells.Item($r, 8)   # column H - rate_cri

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2
    $gVal = $gCell.Value2
    $hVal = $hCell.Value2

    # Swap E<->G and F<->H
    Set-CellText $eCell $gVal
    Set-CellText $fCell $hVal
    Set-CellText $gCell $eVal
    Set-CellText $hCell $fVal
}
